$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93, pushing existing rows 93:147 down to 94:148.
$ws.Rows("93:93").Insert()

# Populate the newly inserted row 93 with the new weekly record.
$ws.Cells.Item(93, 1).Value2 = 8
$ws.Cells.Item(93, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(93, 3).Value2 = "Coquimbo"
$ws.Cells.Item(93, 4).Value2 = 44719
$ws.Cells.Item(93, 5).Value2 = 4
$ws.Cells.Item(93, 6).Value2 = 100112040
$ws.Cells.Item(93, 7).Value2 = "Cilantro"
$ws.Cells.Item(93, 8).Value2 = "Sin especificar"
$ws.Cells.Item(93, 9).Value2 = "Primera"
$ws.Cells.Item(93, 10).Value2 = 3200
$ws.Cells.Item(93, 11).Value2 = 1500
$ws.Cells.Item(93, 12).Value2 = 2000
$ws.Cells.Item(93, 13).Value2 = 1750
$ws.Cells.Item(93, 14).Value2 = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(93, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(93, 16).Value2 = 1167
$ws.Cells.Item(93, 17).Value2 = 1.5
$ws.Cells.Item(93, 18).Value2 = "Hortaliza"
